$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.917.67"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.599.24"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'523.74"
$ws.Range("E5").Value = "  +3.66%  "
$ws.Range("D6").Value = "'154.95"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("D9").Value = "'6.69"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "3.052.13"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "60.918.58"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "'21.72"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "2.600.15"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "'4.76"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'354.90"
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("D20").Value = "'10.58"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "'6.23"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'60.88"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "2.711.66"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "0.0₃0847"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "'7.40"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'6.31"
$ws.Range("E31").Value = "  +10.93%  "
$ws.Range("D32").Value = "'19.39"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").Value = "'147.80"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D35").Value = "'4.18"
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("D36").Value = "'0.935"
$ws.Range("E36").Value = "  +10.50%  "
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.51"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.863"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "'36.48"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "'288.19"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").Value = "'0.623"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "'19.64"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "'4.93"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "'0.0237"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "'19.14"
$ws.Range("E51").Value = "  +8.85%  "
